# Daily update & bug fix
# Appends a new row (row 42) of data to each of the 5 worksheets,
# for date serial 43944 (2020-04-23).

$wb = $excel.ActiveWorkbook

$data = @{
    1 = @(18738, 2021, 194, 394, 172, 4177)
    2 = @(13873.80077473827, 10068.36023028597, 1803.356436746158, 3413.10114543328, 1554.546787248884, 12215.50549762086)
    3 = @(638, 84, 7, 10, 23, 123)
    4 = @(472.3815185336223, 418.4771199129248, 65.06956215063455, 86.62693262521016, 207.8754424809554, 359.7096423766736)
    5 = @(484.8204049150721, 508.1507884656942, 55.77391041482961, 83.16185532020178, 148.2242285516377, 336.8987870064455)
}

for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = 42

    # Column A: date serial, same style as the cells above it (A41)
    $aCell = $ws.Cells.Item($row, 1)
    $ws.Cells.Item($row - 1, 1).Copy($aCell)
    $aCell.Value = 43944

    $vals = $data[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $vals[$c]
    }
}
